$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the style of the existing header cell (E1) onto the new header cells
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Boolean outlier flag values for rows 2-23, columns F (KNN), G (SVM), H (RF)
$values = @{
    2  = @(0,0,0)
    3  = @(0,0,0)
    4  = @(0,0,0)
    5  = @(0,0,0)
    6  = @(0,0,0)
    7  = @(0,0,0)
    8  = @(0,0,0)
    9  = @(1,1,1)
    10 = @(1,1,1)
    11 = @(0,0,0)
    12 = @(0,0,0)
    13 = @(0,0,0)
    14 = @(0,1,0)
    15 = @(0,0,0)
    16 = @(0,0,0)
    17 = @(0,0,0)
    18 = @(1,1,1)
    19 = @(0,0,0)
    20 = @(0,1,0)
    21 = @(0,0,0)
    22 = @(0,0,0)
    23 = @(1,1,1)
}

foreach ($row in $values.Keys) {
    $vals = $values[$row]
    $ws.Cells.Item($row, 6).Value = [bool]($vals[0])
    $ws.Cells.Item($row, 7).Value = [bool]($vals[1])
    $ws.Cells.Item($row, 8).Value = [bool]($vals[2])
}
